# Project Report Part 4 - Investigations: "B. Investigation Design" section added
# plus small wording tweaks to the existing "A. Idea Generation" paragraphs and a
# page-margin tweak.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Small wording fixes inside the existing "Idea Generation" narrative
#    paragraph (currently Paragraphs(3)).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "discussed, recorded, and written down",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "discussed and written down", 2) | Out-Null

$d.Content.Find.Execute(
    "a performance test and analysis every time",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "a performance test and an analysis every time", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Small wording fixes inside the "first meeting" narrative paragraph
#    (currently Paragraphs(4)).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "hardware, software, user interface, organization, etc. and we deliberately",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "hardware, software, organization, etc. We deliberately", 2) | Out-Null

$d.Content.Find.Execute(
    "was the most efficient way to generate ideas",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "was the most effective way to generate ideas", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of "A. Idea Generation" to the
#    end of "B. Investigation Design" (adding a bookmark with this reserved
#    name automatically replaces any existing one, matching Word's own
#    single-instance behaviour for _GoBack).
# ---------------------------------------------------------------------------
$bHeading = $d.Paragraphs(5)
$bEnd = $bHeading.Range
$bEnd.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bEnd) | Out-Null

# ---------------------------------------------------------------------------
# 4. Turn the trailing empty "tab only" paragraph into the new
#    "Investigation Design" narrative paragraph, then append two more
#    trailing paragraphs (a lone space, then a tab + space) to match the
#    final document shape.
# ---------------------------------------------------------------------------
$designParaIndex = $d.Paragraphs.Count
$r = $d.Paragraphs($designParaIndex).Range
$r.Collapse(0)

$r.InsertAfter("The design of our reflow oven controller was established through careful research")
$r.Collapse(0)
$r.InsertAfter(" of the specific component data sheets ")
$r.Collapse(0)
$r.InsertAfter("and examination of the project files and lecture slides provided in the UBC connect website.")
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("Individually, we gathered information for our own tasks and later discussed our plan with the group. Every member")
$r.Collapse(0)
$r.InsertAfter([char]0x2019 + "s plan of action is recorded into a file that was stored in an online tool called ")
$r.Collapse(0)
$r.InsertAfter("GitHub")
$r.Font.Italic = $true
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Font.Italic = $true
$r.Collapse(0)
$r.InsertAfter("(please read 2.3 for further details), which allows access to all the group members for reference. Once the group has heard and approved of the proposal, we continued to work on the design until completion. As aforementioned, each design underwent a performance test and an analysis whenever it is completed in order to keep track of the design's function and operation. During the performance test and analysis, each member of our group was present to observe and allow individual suggestions and constructive criticism for the design to improve on. Again, these are all recorded and stored into GitHub.  ")

# New trailing paragraph #1: a single space.
$d.Paragraphs($designParaIndex).Range.InsertParagraphAfter()
$spaceParaIndex = $designParaIndex + 1
$sr = $d.Paragraphs($spaceParaIndex).Range
$sr.Collapse(0)
$sr.InsertAfter(" ")

# New trailing paragraph #2: a tab followed by a space.
$d.Paragraphs($spaceParaIndex).Range.InsertParagraphAfter()
$tabParaIndex = $spaceParaIndex + 1
$tr = $d.Paragraphs($tabParaIndex).Range
$tr.Collapse(0)
$tr.InsertAfter([char]9)
$tr.Collapse(0)
$tr.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 5. Page margin tweak: bottom margin 1276 -> 1134 twips (63.8pt -> 56.7pt).
# ---------------------------------------------------------------------------
$d.Sections(1).PageSetup.BottomMargin = 56.7
